$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Weekly Quantity"
# Remove the rows for order-weeks 2023-06-11, 2023-06-18, 2023-07-09,
# 2023-07-16, 2023-07-23, 2023-09-30 and 2023-10-07 (penalty/reward cleanup),
# and correct the requested quantity for 2023-06-25 from 52 to 30.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Weekly Quantity")

# Delete rows bottom-to-top so row indices of rows not yet processed stay valid.
$ws1.Rows.Item(24).EntireRow.Delete() | Out-Null
$ws1.Rows.Item(23).EntireRow.Delete() | Out-Null
$ws1.Rows.Item(16).EntireRow.Delete() | Out-Null
$ws1.Rows.Item(15).EntireRow.Delete() | Out-Null
$ws1.Rows.Item(14).EntireRow.Delete() | Out-Null
$ws1.Rows.Item(10).EntireRow.Delete() | Out-Null
$ws1.Rows.Item(9).EntireRow.Delete() | Out-Null

# After the deletions above, the row that held 45102.99999999999 is now row 9.
$ws1.Cells.Item(9, 2).Value = 30

# ---------------------------------------------------------------------------
# Sheet 2: "Monthly Trend"
# Remove the row for order-month 2024-03 (value 72), and correct the
# requested quantities for 2023-06 (450 -> 66) and 2023-07 (836 -> 2).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Monthly Trend")

$ws2.Cells.Item(4, 2).Value = 66
$ws2.Cells.Item(5, 2).Value = 2
$ws2.Rows.Item(10).EntireRow.Delete() | Out-Null
